$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: fix typo "бОродE" -> "бОроду" (B4 "борОду" unchanged)
$ws.Range("A4").Value = "бОроду"

# Row 11: replace "хелло"/"бай" with "лОктя"/"локтЯ"
$ws.Range("A11").Value = "лОктя"
$ws.Range("B11").Value = "локтЯ"

# New row 12: "свЕкла" / "свеклА"
$ws.Range("A12").Value = "свЕкла"
$ws.Range("B12").Value = "свеклА"

# Update selection to match target (B12)
$ws.Range("B12").Select()
